$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 573, shifting existing rows 573:659 down to 574:660
$ws.Rows.Item(573).Insert()

# Populate the newly inserted row 573 with the new record's data
$ws.Cells.Item(573, 1).Value  = 10
$ws.Cells.Item(573, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(573, 3).Value  = "La Araucanía"
$ws.Cells.Item(573, 4).Value  = 45218
$ws.Cells.Item(573, 5).Value  = 9
$ws.Cells.Item(573, 6).Value  = "Fruta"
$ws.Cells.Item(573, 7).Value  = 100108
$ws.Cells.Item(573, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(573, 9).Value  = 100108002
$ws.Cells.Item(573, 10).Value = "Mango"
$ws.Cells.Item(573, 11).Value = "Sin especificar"
$ws.Cells.Item(573, 12).Value = "Primera"
$ws.Cells.Item(573, 13).Value = 2000
$ws.Cells.Item(573, 14).Value = 9000
$ws.Cells.Item(573, 15).Value = 10000
$ws.Cells.Item(573, 16).Value = 9750
$ws.Cells.Item(573, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(573, 18).Value = "Brasil"
$ws.Cells.Item(573, 19).Value = 2438
$ws.Cells.Item(573, 20).Value = 4
